$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 743.25
$ws.Range("I2").Value = 690.1818
$ws.Range("J2").Value = 1327
$ws.Range("K2").Value = 690.1818
$ws.Range("L2").Value = 1327
$ws.Range("M2").Value = -577.1818
$ws.Range("N2").Value = -1553
$ws.Range("H32").Value = 7679.8
$ws.Range("I32").Value = 7600
$ws.Range("J32").Value = 7799.5
$ws.Range("K32").Value = 7600
$ws.Range("L32").Value = 7799.5
$ws.Range("M32").Value = -7274
$ws.Range("N32").Value = -8451.5
$ws.Range("H64").Value = 42668
$ws.Range("I64").Value = 62406.855
$ws.Range("J64").Value = 8125
$ws.Range("K64").Value = 62406.855
$ws.Range("L64").Value = 8125
$ws.Range("M64").Value = -62158.855
$ws.Range("N64").Value = -8621
$ws.Range("H67").Value = 42668
$ws.Range("I67").Value = 62406.855
$ws.Range("J67").Value = 8125
$ws.Range("K67").Value = 62406.855
$ws.Range("L67").Value = 8125
$ws.Range("M67").Value = -61548.855
$ws.Range("N67").Value = -9841
$ws.Range("H80").Value = 118312.35
$ws.Range("I80").Value = 286426.16
$ws.Range("J80").Value = 632.7
$ws.Range("K80").Value = 859278.48
$ws.Range("L80").Value = 1898.1
$ws.Range("M80").Value = -858280.48
$ws.Range("N80").Value = -3894.1
$ws.Range("H83").Value = 118312.35
$ws.Range("I83").Value = 286426.16
$ws.Range("J83").Value = 632.7
$ws.Range("K83").Value = 2577835.44
$ws.Range("L83").Value = 5694.3
$ws.Range("M83").Value = -2572843.44
$ws.Range("N83").Value = -15678.3
$ws.Range("H115").Value = 2630.889
$ws.Range("I115").Value = 409.75
$ws.Range("J115").Value = 20400
$ws.Range("K115").Value = 1229.25
$ws.Range("L115").Value = 61200
$ws.Range("M115").Value = 337.75
$ws.Range("H116").Value = 4507.143
$ws.Range("I116").Value = 3550
$ws.Range("J116").Value = 4890
$ws.Range("K116").Value = 3550
$ws.Range("L116").Value = 4890
$ws.Range("M116").Value = -108
$ws.Range("N116").Value = -11774
$ws.Range("H118").Value = 1005.8
$ws.Range("I118").Value = 507.25
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 1521.75
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 135.25
$ws.Range("H129").Value = 1247.25
$ws.Range("I129").Value = 1247.25
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3741.75
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1258.25
$ws.Range("H132").Value = 5480.161
$ws.Range("I132").Value = 4920.4585
$ws.Range("J132").Value = 7399.143
$ws.Range("K132").Value = 14761.3755
$ws.Range("L132").Value = 22197.429
$ws.Range("M132").Value = -12231.3755
$ws.Range("N132").Value = -27257.429
$ws.Range("H135").Value = 6541
$ws.Range("I135").Value = 7556.1665
$ws.Range("J135").Value = 450
$ws.Range("K135").Value = 68005.4985
$ws.Range("L135").Value = 4050
$ws.Range("M135").Value = -65470.4985
$ws.Range("H137").Value = 5411.2114
$ws.Range("I137").Value = 6121.25
$ws.Range("J137").Value = 1506
$ws.Range("K137").Value = 18363.75
$ws.Range("L137").Value = 4518
$ws.Range("M137").Value = -15813.75
$ws.Range("H138").Value = 3329.96
$ws.Range("I138").Value = 1354.4445
$ws.Range("J138").Value = 4441.1875
$ws.Range("K138").Value = 4063.3335
$ws.Range("L138").Value = 13323.5625
$ws.Range("M138").Value = 1076.6665
$ws.Range("N138").Value = -23603.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 616
$ws.Range("I4").Value = 273
$ws.Range("J4").Value = 787.5
$ws.Range("K4").Value = 273
$ws.Range("L4").Value = 787.5
$ws.Range("M4").Value = -157
$ws.Range("N4").Value = -1019.5
$ws.Range("H45").Value = 6544.154
$ws.Range("I45").Value = 5384.375
$ws.Range("J45").Value = 8399.799999999999
$ws.Range("K45").Value = 5384.375
$ws.Range("L45").Value = 8399.799999999999
$ws.Range("M45").Value = -5007.375
$ws.Range("H110").Value = 10050.071
$ws.Range("I110").Value = 13466.777
$ws.Range("J110").Value = 3900
$ws.Range("K110").Value = 13466.777
$ws.Range("L110").Value = 3900
$ws.Range("M110").Value = -11421.777
$ws.Range("N110").Value = -7990
$ws.Range("H122").Value = 345566.88
$ws.Range("I122").Value = 2349.7778
$ws.Range("J122").Value = 1007485.56
$ws.Range("K122").Value = 7049.3334
$ws.Range("L122").Value = 3022456.68
$ws.Range("M122").Value = -4599.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6002.625
$ws.Range("I94").Value = 5860.143
$ws.Range("J94").Value = 7000
$ws.Range("K94").Value = 5860.143
$ws.Range("L94").Value = 7000
$ws.Range("M94").Value = -5409.143
$ws.Range("H134").Value = 2660.4138
$ws.Range("I134").Value = 2159.923
$ws.Range("J134").Value = 6998
$ws.Range("K134").Value = 6479.768999999999
$ws.Range("L134").Value = 20994
$ws.Range("M134").Value = -3944.768999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1326.25
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1326.25
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1326.25
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2026.25
$ws.Range("H31").Value = 1816.5834
$ws.Range("I31").Value = 1254.4546
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 1254.4546
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -959.4546
$ws.Range("N31").Value = -8590
$ws.Range("H32").Value = 674.1875
$ws.Range("I32").Value = 682.25
$ws.Range("J32").Value = 650
$ws.Range("K32").Value = 682.25
$ws.Range("L32").Value = 650
$ws.Range("M32").Value = -366.25
$ws.Range("N32").Value = -1282
$ws.Range("H34").Value = 1816.5834
$ws.Range("I34").Value = 1254.4546
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 1254.4546
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -1052.4546
$ws.Range("N34").Value = -8404
$ws.Range("H35").Value = 3576.5
$ws.Range("I35").Value = 4953
$ws.Range("J35").Value = 2200
$ws.Range("K35").Value = 4953
$ws.Range("L35").Value = 2200
$ws.Range("M35").Value = -4659
$ws.Range("N35").Value = -2788
$ws.Range("H58").Value = 4160.75
$ws.Range("I58").Value = 1594.75
$ws.Range("J58").Value = 5016.0835
$ws.Range("K58").Value = 1594.75
$ws.Range("L58").Value = 5016.0835
$ws.Range("M58").Value = -1391.75
$ws.Range("N58").Value = -5422.0835
$ws.Range("H99").Value = 7267669
$ws.Range("I99").Value = 14530213
$ws.Range("J99").Value = 5124.875
$ws.Range("K99").Value = 14530213
$ws.Range("L99").Value = 5124.875
$ws.Range("M99").Value = -14528715
$ws.Range("N99").Value = -8120.875
$ws.Range("H126").Value = 7267669
$ws.Range("I126").Value = 14530213
$ws.Range("J126").Value = 5124.875
$ws.Range("K126").Value = 43590639
$ws.Range("L126").Value = 15374.625
$ws.Range("M126").Value = -43588169
$ws.Range("N126").Value = -20314.625
$ws.Range("H132").Value = 28736.904
$ws.Range("I132").Value = 9153.5625
$ws.Range("J132").Value = 91403.60000000001
$ws.Range("K132").Value = 27460.6875
$ws.Range("L132").Value = 274210.8
$ws.Range("M132").Value = -24930.6875
$ws.Range("N132").Value = -279270.8
$ws.Range("H136").Value = 4160.75
$ws.Range("I136").Value = 1594.75
$ws.Range("J136").Value = 5016.0835
$ws.Range("K136").Value = 4784.25
$ws.Range("L136").Value = 15048.2505
$ws.Range("M136").Value = -2234.25
$ws.Range("N136").Value = -20148.2505
$ws.Range("H138").Value = 145000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 145000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 145000
$ws.Range("N138").Value = -155280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 124117.8
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 124117.8
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 372353.4
$ws.Range("N37").Value = -372577.4
$ws.Range("H68").Value = 10222.667
$ws.Range("I68").Value = 2738.3333
$ws.Range("J68").Value = 12093.75
$ws.Range("K68").Value = 8214.999899999999
$ws.Range("L68").Value = 36281.25
$ws.Range("M68").Value = -7403.999899999999
$ws.Range("N68").Value = -37903.25
$ws.Range("H71").Value = 10222.667
$ws.Range("I71").Value = 2738.3333
$ws.Range("J71").Value = 12093.75
$ws.Range("K71").Value = 24644.9997
$ws.Range("L71").Value = 108843.75
$ws.Range("M71").Value = -20588.9997
$ws.Range("N71").Value = -116955.75
$ws.Range("H80").Value = 119392
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 119392
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 358176
$ws.Range("N80").Value = -360048
$ws.Range("H83").Value = 119392
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 119392
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 1074528
$ws.Range("N83").Value = -1083888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6953.625
$ws.Range("I70").Value = 6308
$ws.Range("J70").Value = 7455.778
$ws.Range("K70").Value = 6308
$ws.Range("L70").Value = 7455.778
$ws.Range("M70").Value = -6038
$ws.Range("N70").Value = -7995.778
$ws.Range("H73").Value = 6953.625
$ws.Range("I73").Value = 6308
$ws.Range("J73").Value = 7455.778
$ws.Range("K73").Value = 6308
$ws.Range("L73").Value = 7455.778
$ws.Range("M73").Value = -5372
$ws.Range("N73").Value = -9327.778
$ws.Range("H122").Value = 10898.68
$ws.Range("I122").Value = 9261.474
$ws.Range("J122").Value = 16083.167
$ws.Range("K122").Value = 27784.422
$ws.Range("L122").Value = 48249.501
$ws.Range("M122").Value = -25334.422
$ws.Range("N122").Value = -53149.501
$ws.Range("H132").Value = 3487.2173
$ws.Range("I132").Value = 3513.9546
$ws.Range("J132").Value = 2899
$ws.Range("K132").Value = 10541.8638
$ws.Range("L132").Value = 8697
$ws.Range("M132").Value = -8011.863799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3740.1667
$ws.Range("I32").Value = 3740.1667
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3740.1667
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3423.1667
$ws.Range("H136").Value = 10838.95
$ws.Range("I136").Value = 13778.7
$ws.Range("J136").Value = 7899.2
$ws.Range("K136").Value = 41336.10000000001
$ws.Range("L136").Value = 23697.6
$ws.Range("M136").Value = -38786.10000000001
$ws.Range("N136").Value = -28797.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H75").Value = 22250
$ws.Range("I75").Value = 19666.666
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 19666.666
$ws.Range("L75").Value = 30000
$ws.Range("M75").Value = -18730.666
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 22250
$ws.Range("I78").Value = 19666.666
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 58999.99800000001
$ws.Range("L78").Value = 90000
$ws.Range("M78").Value = -54319.99800000001
$ws.Range("N78").Value = -99360
$ws.Range("H96").Value = 4288.4614
$ws.Range("I96").Value = 4137.625
$ws.Range("J96").Value = 4529.8
$ws.Range("K96").Value = 4137.625
$ws.Range("L96").Value = 4529.8
$ws.Range("M96").Value = -2764.625
$ws.Range("N96").Value = -7275.8
$ws.Range("H126").Value = 17574.963
$ws.Range("I126").Value = 21536.95
$ws.Range("J126").Value = 6255
$ws.Range("K126").Value = 64610.85000000001
$ws.Range("L126").Value = 18765
$ws.Range("M126").Value = -62140.85000000001
$ws.Range("N126").Value = -23705
$ws.Range("H136").Value = 2321.1875
$ws.Range("I136").Value = 1753.25
$ws.Range("J136").Value = 4025
$ws.Range("K136").Value = 5259.75
$ws.Range("L136").Value = 12075
$ws.Range("M136").Value = -2709.75
$ws.Range("N136").Value = -17175


Write-Host "Applied all changes"
